$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 9 - update kickoff time text
$ws.Range("C9").Value = "17:30"

# Row 10 - updated odds
$ws.Range("G10").Value = 1.73
$ws.Range("H10").Value = 3.9
$ws.Range("I10").Value = 4.2
$ws.Range("J10").Value = 1.04
$ws.Range("K10").Value = 9.5
$ws.Range("L10").Value = 1.2
$ws.Range("M10").Value = 4.25
$ws.Range("N10").Value = 1.67
$ws.Range("O10").Value = 2.15
$ws.Range("P10").Value = 1.32
$ws.Range("Q10").Value = 3.2
$ws.Range("R10").Value = 1.57
$ws.Range("S10").Value = 2.25
$ws.Range("T10").Value = 8.5
$ws.Range("U10").Value = 9
$ws.Range("AB10").Value = 15
$ws.Range("AE10").Value = 15
$ws.Range("AF10").Value = 23
$ws.Range("AI10").Value = 34
$ws.Range("AJ10").Value = 34

# Row 18 - updated odds
$ws.Range("G18").Value = 1.5
$ws.Range("I18").Value = 5
$ws.Range("N18").Value = 1.5
$ws.Range("O18").Value = 2.5
$ws.Range("T18").Value = 9.5
$ws.Range("V18").Value = 9
$ws.Range("AA18").Value = 9.5
$ws.Range("AE18").Value = 19
$ws.Range("AF18").Value = 29
$ws.Range("AJ18").Value = 34

$wb.Save()
